$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that look like plain numbers need a leading apostrophe (quote-prefix)
# so Excel stores them as text, matching the original inlineStr cell type,
# instead of silently converting them to numeric values.

$ws.Range("D2").Value = "63.472.38"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.632.77"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'606.22"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").Value = "'147.24"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("D10").Value = "'5.58"
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").Value = "'0.372"
$ws.Range("E11").Value = "  +4.79%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "'27.55"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "3.109.20"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "63.324.14"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("D17").Value = "2.656.60"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'11.56"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("E19").Value = "  +5.26%  "
$ws.Range("D20").Value = "'344.92"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "'6.88"
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").Value = "'66.94"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").Value = "'9.07"
$ws.Range("E26").Value = "  +7.81%  "
$ws.Range("D27").Value = "'1.57"
$ws.Range("E27").Value = "  +2.97%  "
$ws.Range("D28").Value = "'565.41"
$ws.Range("E28").Value = "  +7.38%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'8.06"
$ws.Range("E29").Value = "  +3.63%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.163"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("D33").Value = "0.0₃0851"
$ws.Range("E33").Value = "  +5.94%  "
$ws.Range("D34").Value = "'1.76"
$ws.Range("E34").Value = "  -2.86%  "
$ws.Range("D35").Value = "'5.18"
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("D36").Value = "'167.43"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "'1.96"
$ws.Range("E39").Value = "  +9.55%  "
$ws.Range("D40").Value = "'19.15"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'166.21"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("D44").Value = "'22.20"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  +3.97%  "
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").Value = "'1.93"
$ws.Range("E49").Value = "  +13.72%  "
$ws.Range("D50").Value = "'18.88"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("E51").Value = "  +6.36%  "
